$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $name = $ws.Cells.Item($r, 1).Value2
    if ([string]::IsNullOrEmpty($name)) { continue }

    $rest = $name -replace '^Rule - ', ''
    $firstWord = ($rest -split ' ')[0]
    if ($firstWord -like '*-*') {
        $firstWord = ($firstWord -split '-')[0]
    }
    $normalized = "rule_" + $firstWord.ToLower()

    $ws.Cells.Item($r, 2).Value2 = $normalized
}
